$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("G2").Value = 4.2
$ws.Range("I2").Value = 2
$ws.Range("L2").Value = 2.75
$ws.Range("AH2").Value = 51

# Row 3 updates
$ws.Range("G3").Value = 1.85
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2.6
$ws.Range("K3").Value = 1.95
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.93
$ws.Range("U3").Value = 4.1
$ws.Range("V3").Value = 1.24
$ws.Range("AD3").Value = 7.5
$ws.Range("AF3").Value = 15
$ws.Range("AG3").Value = 19

# Row 5 updates
$ws.Range("G5").Value = 2.63
$ws.Range("I5").Value = 2.75
$ws.Range("J5").Value = 3.4
$ws.Range("K5").Value = 2.05
$ws.Range("W5").Value = 4
$ws.Range("X5").Value = 1.25
$ws.Range("AI5").Value = 8.5
$ws.Range("AM5").Value = 301

# Row 6 updates
$ws.Range("G6").Value = 3.2
$ws.Range("H6").Value = 3.2
$ws.Range("J6").Value = 3.75
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.25
$ws.Range("S6").Value = 2.05
$ws.Range("T6").Value = 1.75
$ws.Range("W6").Value = 3.5
$ws.Range("X6").Value = 1.29
$ws.Range("Y6").Value = 1.44
$ws.Range("Z6").Value = 2.63
$ws.Range("AA6").Value = 1.8
$ws.Range("AB6").Value = 1.91
$ws.Range("AC6").Value = 9.5
$ws.Range("AG6").Value = 26
$ws.Range("AH6").Value = 34
$ws.Range("AI6").Value = 9
$ws.Range("AM6").Value = 251
$ws.Range("AN6").Value = 7.5
$ws.Range("AO6").Value = 11
$ws.Range("AP6").Value = 9.5
$ws.Range("AR6").Value = 19
$ws.Range("AS6").Value = 29

# Row 7 updates
$ws.Range("G7").Value = 3.2
$ws.Range("I7").Value = 2.2
$ws.Range("S7").Value = 1.93
$ws.Range("T7").Value = 1.93
$ws.Range("W7").Value = 3.25
$ws.Range("X7").Value = 1.33
$ws.Range("AK7").Value = 13
$ws.Range("AL7").Value = 41
$ws.Range("AO7").Value = 11
$ws.Range("AQ7").Value = 21

